# Updates cryptos list: price (D) and volume-change (E) columns
# Cells must remain plain text to match the original inline-string formatting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "58.526.81"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -2.72%  "
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.289.49"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -5.42%  "
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "546.43"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -1.30%  "
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "130.92"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -4.45%  "
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.574"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -2.37%  "
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.286.87"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -5.41%  "
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.102"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -3.43%  "
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.54"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -2.23%  "
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +0.26%  "
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.335"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -5.21%  "
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "23.64"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -5.06%  "
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.694.21"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -5.54%  "
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "58.496.70"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -2.63%  "
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -3.30%  "
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.287.60"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -6.30%  "
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "10.65"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -5.66%  "
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.30"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -3.95%  "
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "315.28"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -3.94%  "
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -4.72%  "
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "62.81"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -3.88%  "
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.171"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -3.00%  "
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -0.57%  "
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.14"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -6.85%  "
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.31"
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -5.65%  "
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.75"
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.72%  "
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "170.67"
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0723"
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -6.30%  "
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.ClearFormats()
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.78"
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -5.23%  "
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -4.95%  "
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "17.84"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -3.78%  "
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -5.14%  "
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.94"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -6.09%  "
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "37.94"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -1.90%  "
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -5.32%  "
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "299.83"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -8.88%  "
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "140.99"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.97%  "
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.45"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -5.68%  "
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0949"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -1.81%  "
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -3.41%  "
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.556"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -3.30%  "
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "18.43"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -7.78%  "
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -3.92%  "
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "16.65"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -5.04%  "
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "
$cell.ClearFormats()
